# SnapperCodingChallenge V1.1 To Do List - folder/content rename & reshuffle
# (commit: "Fixed folder names to reflect csproj titles.")

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# -----------------------------------------------------------------
# Sheet2 ("To Do" detail sheet): move the existing three tasks down
# into their new rows, then fill in the newly added task list plus
# the new "Actors" mini-table in column R.
# -----------------------------------------------------------------

# Preserve the text of the three existing tasks before we overwrite them
$setOptionsText   = $ws2.Range("B16").Value2
$loadSnapperText  = $ws2.Range("B21").Value2
$loadTargetText   = $ws2.Range("B25").Value2

# Move them to their new homes
$ws2.Range("B6").Value  = $setOptionsText
$ws2.Range("B8").Value  = $loadSnapperText
$ws2.Range("B10").Value = $loadTargetText

# Clear out the now-vacated old cells
$ws2.Range("B21").ClearContents()
$ws2.Range("B25").ClearContents()

# New tasks appended to the bottom of the list
$ws2.Range("B12").Value = "Pause the Console to Allow User to Eamine Data "
$ws2.Range("B14").Value = "Summarise TargetsFound "
$ws2.Range("B16").Value = "Write Output File "
$ws2.Range("B18").Value = "Write "

# New "Actors" mini table in column R
$ws2.Range("R3").Value = "Actors "
$ws2.Range("R3").Font.Bold = $true
$ws2.Range("R5").Value = "UI - Something to echo whats actually going on "
$ws2.Range("R6").Value = "Dump - Something to dump whats being echoed to "

# -----------------------------------------------------------------
# Reposition the two logo pictures on Sheet2 (moved up/right)
# -----------------------------------------------------------------
$pic1 = $ws2.Shapes.Item(1)
$pic1.Left   = 332.25
$pic1.Top    = 15.75
$pic1.Width  = 428.1964566929134
$pic1.Height = 166.4792125984252

$pic2 = $ws2.Shapes.Item(2)
$pic2.Left   = 330.75
$pic2.Top    = 186.75007874015748
$pic2.Width  = 429.75
$pic2.Height = 54.33622047244094

# -----------------------------------------------------------------
# View state: Sheet2 selection moves to R7; Sheet1 becomes the
# active/selected tab with its selection moved to AC44.
# -----------------------------------------------------------------
$ws2.Activate() | Out-Null
$ws2.Range("R7").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("AC44").Select() | Out-Null
